$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Active Report")
$ws.Range("C2").Value = "00 46 09 ? * MON,TUE,WED,THU,FRI"
